# Scheduled-runner refresh: updates computed market-price / profit
# columns (H, I, J, K, L, M, N) across the Leve-profit sheets.
# Values below mirror the latest Universalis price pull for each
# sheet's tracked items/recipes.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4065
$ws.Range("I106").Value = 4312.5
$ws.Range("K106").Value = 4312.5
$ws.Range("M106").Value = -3681.5
$ws.Range("H116").Value = 4397.625
$ws.Range("I116").Value = 3631.8333
$ws.Range("J116").Value = 6695
$ws.Range("K116").Value = 3631.8333
$ws.Range("L116").Value = 6695
$ws.Range("M116").Value = -189.8332999999998
$ws.Range("N116").Value = -13579
$ws.Range("H132").Value = 11179.066
$ws.Range("I132").Value = 8682.406999999999
$ws.Range("K132").Value = 26047.221
$ws.Range("M132").Value = -23517.221
$ws.Range("H133").Value = 99900
$ws.Range("J133").Value = 99900
$ws.Range("L133").Value = 99900
$ws.Range("N133").Value = -110020

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3297
$ws.Range("I45").Value = 3297
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3297
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2920
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 75089.66
$ws.Range("I61").Value = 3401.64
$ws.Range("K61").Value = 3401.64
$ws.Range("M61").Value = -3189.64
$ws.Range("H122").Value = 5650.905
$ws.Range("I122").Value = 4822
$ws.Range("K122").Value = 14466
$ws.Range("M122").Value = -12016
$ws.Range("H136").Value = 75089.66
$ws.Range("I136").Value = 3401.64
$ws.Range("K136").Value = 10204.92
$ws.Range("M136").Value = -7654.92

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1009.3871
$ws.Range("I94").Value = 1085.2693
$ws.Range("K94").Value = 1085.2693
$ws.Range("M94").Value = -634.2692999999999
$ws.Range("H138").Value = 85000
$ws.Range("J138").Value = 85000
$ws.Range("L138").Value = 85000
$ws.Range("N138").Value = -95280

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 22623
$ws.Range("J59").Value = 29826.875
$ws.Range("L59").Value = 29826.875
$ws.Range("N59").Value = -32116.875
$ws.Range("H86").Value = 10484.591
$ws.Range("J86").Value = 7165.6665
$ws.Range("L86").Value = 7165.6665
$ws.Range("N86").Value = -9411.666499999999
$ws.Range("H89").Value = 10484.591
$ws.Range("J89").Value = 7165.6665
$ws.Range("L89").Value = 35828.3325
$ws.Range("N89").Value = -47060.3325
$ws.Range("H99").Value = 4256.143
$ws.Range("I99").Value = 4132.1665
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 4132.1665
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -2634.1665
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 4256.143
$ws.Range("I126").Value = 4132.1665
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 12396.4995
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -9926.499500000002
$ws.Range("N126").Value = -19940
$ws.Range("H134").Value = 28577200
$ws.Range("I134").Value = 2113.4119
$ws.Range("K134").Value = 6340.2357
$ws.Range("M134").Value = -3805.2357

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 321.2143
$ws.Range("J35").Value = 300
$ws.Range("L35").Value = 900
$ws.Range("N35").Value = -1476
$ws.Range("H68").Value = 1608.2222
$ws.Range("I68").Value = 2347.6667
$ws.Range("J68").Value = 1238.5
$ws.Range("K68").Value = 7043.000100000001
$ws.Range("L68").Value = 3715.5
$ws.Range("M68").Value = -6232.000100000001
$ws.Range("N68").Value = -5337.5
$ws.Range("H71").Value = 1608.2222
$ws.Range("I71").Value = 2347.6667
$ws.Range("J71").Value = 1238.5
$ws.Range("K71").Value = 21129.0003
$ws.Range("L71").Value = 11146.5
$ws.Range("M71").Value = -17073.0003
$ws.Range("N71").Value = -19258.5
$ws.Range("H75").Value = 1366.6666
$ws.Range("I75").Value = 600
$ws.Range("K75").Value = 1800
$ws.Range("M75").Value = -802
$ws.Range("H78").Value = 1366.6666
$ws.Range("I78").Value = 600
$ws.Range("K78").Value = 5400
$ws.Range("M78").Value = -408

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11497.5
$ws.Range("I80").Value = 14990
$ws.Range("J80").Value = 10333.333
$ws.Range("K80").Value = 14990
$ws.Range("L80").Value = 10333.333
$ws.Range("M80").Value = -13992
$ws.Range("N80").Value = -12329.333
$ws.Range("H83").Value = 11497.5
$ws.Range("I83").Value = 14990
$ws.Range("J83").Value = 10333.333
$ws.Range("K83").Value = 74950
$ws.Range("L83").Value = 51666.665
$ws.Range("M83").Value = -69958
$ws.Range("N83").Value = -61650.665

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4532.3335
$ws.Range("I46").Value = 4999
$ws.Range("J46").Value = 4460.5386
$ws.Range("K46").Value = 4999
$ws.Range("L46").Value = 4460.5386
$ws.Range("M46").Value = -4811
$ws.Range("N46").Value = -4836.5386
$ws.Range("H68").Value = 4411
$ws.Range("I68").Value = 4411
$ws.Range("K68").Value = 4411
$ws.Range("M68").Value = -3662
$ws.Range("H71").Value = 4411
$ws.Range("I71").Value = 4411
$ws.Range("K71").Value = 22055
$ws.Range("M71").Value = -18311
$ws.Range("H132").Value = 1915120.1
$ws.Range("I132").Value = 5056.875
$ws.Range("K132").Value = 15170.625
$ws.Range("M132").Value = -12640.625

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5517500
$ws.Range("J3").Value = 35000
$ws.Range("L3").Value = 35000
$ws.Range("N3").Value = -35228
$ws.Range("H62").Value = 19815.154
$ws.Range("I62").Value = 44999.5
$ws.Range("J62").Value = 15236.182
$ws.Range("K62").Value = 44999.5
$ws.Range("L62").Value = 15236.182
$ws.Range("M62").Value = -44375.5
$ws.Range("N62").Value = -16484.182
$ws.Range("H65").Value = 19815.154
$ws.Range("I65").Value = 44999.5
$ws.Range("J65").Value = 15236.182
$ws.Range("K65").Value = 224997.5
$ws.Range("L65").Value = 76180.91
$ws.Range("M65").Value = -221877.5
$ws.Range("N65").Value = -82420.91
$ws.Range("H132").Value = 9839.352999999999
$ws.Range("I132").Value = 4366.8623
$ws.Range("K132").Value = 13100.5869
$ws.Range("M132").Value = -10570.5869
$ws.Range("H136").Value = 13012.9
$ws.Range("I136").Value = 1590.88
$ws.Range("K136").Value = 4772.64
$ws.Range("M136").Value = -2222.64
